$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B41: change from text "2" to numeric 2
$ws.Range("B41").Value = 2

# Add new row 42 with data
$ws.Range("A42").Value = "Sunsi Wu"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "3"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "无"
$ws.Range("D42").Value = "CRT"
$ws.Range("E42").Value = "MET"
$ws.Range("F42").Value = "cf97de89-8b46-4ca2-a071-801296a106cf"
$ws.Range("G42").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H42").Value = "Some important implementation details are missing (activation functions, loss function used), and others have to be deduced by observing the output dimensions of the individual layers of the network."
